# cambios de las fracciones
# Update the reporting quarter on the main "Reporte de Formatos" sheet:
#   B8 (fecha inicio)  : 2022-07-01 -> 2022-10-01
#   C8 (fecha termino) : 2022-09-30 -> 2022-12-31
#   H8 (fecha validacion)   : 2022-10-10 -> 2023-01-10
#   I8 (fecha actualizacion): 2022-10-10 -> 2023-01-10

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Reporte de Formatos")

$ws.Range("B8").Value = (Get-Date -Year 2022 -Month 10 -Day 1 -Hour 0 -Minute 0 -Second 0).Date
$ws.Range("C8").Value = (Get-Date -Year 2022 -Month 12 -Day 31 -Hour 0 -Minute 0 -Second 0).Date
$ws.Range("H8").Value = (Get-Date -Year 2023 -Month 1 -Day 10 -Hour 0 -Minute 0 -Second 0).Date
$ws.Range("I8").Value = (Get-Date -Year 2023 -Month 1 -Day 10 -Hour 0 -Minute 0 -Second 0).Date

# Re-apply the thin-border-only (no special alignment) style to the
# previously-blank J8 cell, matching the reformat that came with the
# quarter change.
$j8 = $ws.Range("J8")
$j8.ClearFormats()
$j8.Borders.LineStyle = 1

# The same border-only style gets re-applied to the first data row of each
# of the three child ("cargo") tables.
$wsTabla1 = $wb.Worksheets.Item("Tabla_397514")
$rngTabla1 = $wsTabla1.Range("A4:E4")
$rngTabla1.ClearFormats()
$rngTabla1.Borders.LineStyle = 1

$wsTabla2 = $wb.Worksheets.Item("Tabla_397515")
$rngTabla2 = $wsTabla2.Range("A4:E4")
$rngTabla2.ClearFormats()
$rngTabla2.Borders.LineStyle = 1

$wsTabla3 = $wb.Worksheets.Item("Tabla_397516")
$rngTabla3 = $wsTabla3.Range("A4:E4")
$rngTabla3.ClearFormats()
$rngTabla3.Borders.LineStyle = 1
